$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.14118879991426
$ws.Range("C2").Value = 4.714971814194675
$ws.Range("D2").Value = 7.60243751228304
$ws.Range("E2").Value = 9.938551248749587
$ws.Range("F2").Value = 37.75385985877457
$ws.Range("I2").Value = 31.01830222044006
$ws.Range("K2").Value = 13.9804170890811
$ws.Range("L2").Value = 10.35233306271476
$ws.Range("N2").Value = 22.18326019142626
$ws.Range("B3").Value = 15.94619127431521
$ws.Range("C3").Value = 4.457157147758701
$ws.Range("D3").Value = 7.609906557120055
$ws.Range("E3").Value = 9.947211762593218
$ws.Range("F3").Value = 37.67890564883053
$ws.Range("I3").Value = 31.04199030422711
$ws.Range("K3").Value = 13.84376505895688
$ws.Range("L3").Value = 10.34271338168336
$ws.Range("N3").Value = 22.23616157705348
$ws.Range("B4").Value = 15.82962789282473
$ws.Range("C4").Value = 4.290122461420101
$ws.Range("D4").Value = 7.61460490854741
$ws.Range("E4").Value = 9.954011870861381
$ws.Range("F4").Value = 37.64140246172978
$ws.Range("I4").Value = 31.06210622571427
$ws.Range("K4").Value = 13.76283001884183
$ws.Range("L4").Value = 10.33875354970141
$ws.Range("N4").Value = 22.27054980033538
$ws.Range("B5").Value = 15.78298076116336
$ws.Range("C5").Value = 4.219881956196033
$ws.Range("D5").Value = 7.61654798267427
$ws.Range("E5").Value = 9.957156025520746
$ws.Range("F5").Value = 37.62826945504848
$ws.Range("I5").Value = 31.07170291591019
$ws.Range("K5").Value = 13.73063072870059
$ws.Range("L5").Value = 10.33763088198513
$ws.Range("N5").Value = 22.28504317742432
$ws.Range("B6").Value = 15.77528822334049
$ws.Range("C6").Value = 4.208088178793432
$ws.Range("D6").Value = 7.616872353647246
$ws.Range("E6").Value = 9.957700647828261
$ws.Range("F6").Value = 37.6262187603636
$ws.Range("I6").Value = 31.0733808954544
$ws.Range("K6").Value = 13.72533232458965
$ws.Range("L6").Value = 10.33747415637652
$ws.Range("N6").Value = 22.28747879053438
$ws.Range("B7").Value = 15.82899526354942
$ws.Range("C7").Value = 4.289183933733787
$ws.Range("D7").Value = 7.61463099803536
$ws.Range("E7").Value = 9.954052763247098
$ws.Range("F7").Value = 37.64121663218366
$ws.Range("I7").Value = 31.06222998696264
$ws.Range("K7").Value = 13.76239255424169
$ws.Range("L7").Value = 10.3387364192484
$ws.Range("N7").Value = 22.2707433194002
$ws.Range("B8").Value = 16.07333145965105
$ws.Range("C8").Value = 4.627903070063746
$ws.Range("D8").Value = 7.604989653914703
$ws.Range("E8").Value = 9.941229876895665
$ws.Range("F8").Value = 37.72625270988142
$ws.Range("I8").Value = 31.02531265243648
$ws.Range("K8").Value = 13.93270490890878
$ws.Range("L8").Value = 10.34861313697552
$ws.Range("N8").Value = 22.20110516954001
$ws.Range("B9").Value = 16.57500271281531
$ws.Range("C9").Value = 5.222029863805265
$ws.Range("D9").Value = 7.586964278749349
$ws.Range("E9").Value = 9.92783364937438
$ws.Range("F9").Value = 37.96018208884473
$ws.Range("I9").Value = 30.99718470839059
$ws.Range("K9").Value = 14.28861576209488
$ws.Range("L9").Value = 10.38335069873541
$ws.Range("N9").Value = 22.07965062623
$ws.Range("B10").Value = 16.95377678489797
$ws.Range("C10").Value = 5.615050713781351
$ws.Range("D10").Value = 7.57424382141223
$ws.Range("E10").Value = 9.925132632703599
$ws.Range("F10").Value = 38.17231739666779
$ws.Range("I10").Value = 31.00356771817148
$ws.Range("K10").Value = 14.56119984134638
$ws.Range("L10").Value = 10.41812218228075
$ws.Range("N10").Value = 21.99959563387434
$ws.Range("B11").Value = 17.12754711216996
$ws.Range("C11").Value = 5.784314008354198
$ws.Range("D11").Value = 7.568567353478108
$ws.Range("E11").Value = 9.925448254422937
$ws.Range("F11").Value = 38.27738100427604
$ws.Range("I11").Value = 31.012348960417
$ws.Range("K11").Value = 14.68711586200622
$ws.Range("L11").Value = 10.43591742168709
$ws.Range("N11").Value = 21.96516308156763
$ws.Range("B12").Value = 17.19349710297033
$ws.Range("C12").Value = 5.84703557384842
$ws.Range("D12").Value = 7.566433428164983
$ws.Range("E12").Value = 9.92578913929867
$ws.Range("F12").Value = 38.31837823184227
$ws.Range("I12").Value = 31.01651875540763
$ws.Range("K12").Value = 14.73503021523607
$ws.Range("L12").Value = 10.44293717139402
$ws.Range("N12").Value = 21.95240939979721
$ws.Range("B13").Value = 17.17928805913396
$ws.Range("C13").Value = 5.833588567973509
$ws.Range("D13").Value = 7.566892315299611
$ws.Range("E13").Value = 9.925705890583968
$ws.Range("F13").Value = 38.30949517181374
$ws.Range("I13").Value = 31.01558317066632
$ws.Range("K13").Value = 14.72470133460022
$ws.Range("L13").Value = 10.44141289712686
$ws.Range("N13").Value = 21.95514345359141
$ws.Range("B14").Value = 17.13297027344441
$ws.Range("C14").Value = 5.789501724624603
$ws.Range("D14").Value = 7.568391482094808
$ws.Range("E14").Value = 9.925471867104793
$ws.Range("F14").Value = 38.28072969931275
$ws.Range("I14").Value = 31.01267509192085
$ws.Range("K14").Value = 14.69105338647529
$ws.Range("L14").Value = 10.43648932736199
$ws.Range("N14").Value = 21.96410811440585
$ws.Range("B15").Value = 17.1046165574535
$ws.Range("C15").Value = 5.762318103348442
$ws.Range("D15").Value = 7.569311794389817
$ws.Range("E15").Value = 9.925357326765633
$ws.Range("F15").Value = 38.26326724006599
$ws.Range("I15").Value = 31.01100376411549
$ws.Range("K15").Value = 14.67047207875387
$ws.Range("L15").Value = 10.43351000086101
$ws.Range("N15").Value = 21.96963635982558
$ws.Range("B16").Value = 16.94244448914634
$ws.Range("C16").Value = 5.603796818454475
$ws.Range("D16").Value = 7.574616990365722
$ws.Range("E16").Value = 9.925143023620702
$ws.Range("F16").Value = 38.16562192713826
$ws.Range("I16").Value = 31.00311207041865
$ws.Range("K16").Value = 14.55300582779628
$ws.Range("L16").Value = 10.4169987401903
$ws.Range("N16").Value = 22.0018858014637
$ws.Range("B17").Value = 16.84328663980047
$ws.Range("C17").Value = 5.504105443785968
$ws.Range("D17").Value = 7.577899616179802
$ws.Range("E17").Value = 9.925406621661715
$ws.Range("F17").Value = 38.10789951980733
$ws.Range("I17").Value = 30.99977585215267
$ws.Range("K17").Value = 14.48140437743189
$ws.Range("L17").Value = 10.40737386360274
$ws.Range("N17").Value = 22.02217798145908
$ws.Range("B18").Value = 16.78639523560432
$ws.Range("C18").Value = 5.445869337558499
$ws.Range("D18").Value = 7.579798074119591
$ws.Range("E18").Value = 9.925703643892476
$ws.Range("F18").Value = 38.07550630683609
$ws.Range("I18").Value = 30.99841024439669
$ws.Range("K18").Value = 14.44040404724441
$ws.Range("L18").Value = 10.40202423827971
$ws.Range("N18").Value = 22.03403634459954
$ws.Range("B19").Value = 16.76715905078289
$ws.Range("C19").Value = 5.425997852852001
$ws.Range("D19").Value = 7.580442648194738
$ws.Range("E19").Value = 9.925829204417903
$ws.Range("F19").Value = 38.06467771530828
$ws.Range("I19").Value = 30.99804291422302
$ws.Range("K19").Value = 14.42655476404174
$ws.Range("L19").Value = 10.4002450479126
$ws.Range("N19").Value = 22.03808347999323
$ws.Range("B20").Value = 16.85382797403922
$ws.Range("C20").Value = 5.514810565374571
$ws.Range("D20").Value = 7.577549102419383
$ws.Range("E20").Value = 9.92536351688803
$ws.Range("F20").Value = 38.11396076154904
$ws.Range("I20").Value = 31.00007373725108
$ws.Range("K20").Value = 14.48900785193402
$ws.Range("L20").Value = 10.40837918456504
$ws.Range("N20").Value = 22.01999850902014
$ws.Range("B21").Value = 17.14657144384786
$ws.Range("C21").Value = 5.802488441879687
$ws.Range("D21").Value = 7.56795071791945
$ws.Range("E21").Value = 9.925534603763335
$ws.Range("F21").Value = 38.28914608769508
$ws.Range("I21").Value = 31.01350635262717
$ws.Range("K21").Value = 14.70093063055841
$ws.Range("L21").Value = 10.43792789841692
$ws.Range("N21").Value = 21.96146723837939
$ws.Range("B22").Value = 17.33872398274546
$ws.Range("C22").Value = 5.98248827064169
$ws.Range("D22").Value = 7.561768626390535
$ws.Range("E22").Value = 9.926936314453968
$ws.Range("F22").Value = 38.41069358841602
$ws.Range("I22").Value = 31.02720724094737
$ws.Range("K22").Value = 14.84077092944432
$ws.Range("L22").Value = 10.45887630047193
$ws.Range("N22").Value = 21.92487581086601
$ws.Range("B23").Value = 17.23611355096818
$ws.Range("C23").Value = 5.887153509513404
$ws.Range("D23").Value = 7.565059865953919
$ws.Range("E23").Value = 9.926070427078212
$ws.Range("F23").Value = 38.34518282887188
$ws.Range("I23").Value = 31.01944481412897
$ws.Range("K23").Value = 14.76602736731545
$ws.Range("L23").Value = 10.44754715623623
$ws.Range("N23").Value = 21.94425333293485
$ws.Range("B24").Value = 16.84906187134834
$ws.Range("C24").Value = 5.509973652048688
$ws.Range("D24").Value = 7.577707534656806
$ws.Range("E24").Value = 9.925382551370701
$ws.Range("F24").Value = 38.11121800503199
$ws.Range("I24").Value = 30.99993734246425
$ws.Range("K24").Value = 14.48556980605152
$ws.Range("L24").Value = 10.40792410605782
$ws.Range("N24").Value = 22.02098324997936
$ws.Range("B25").Value = 16.43723882974686
$ws.Range("C25").Value = 5.068898992790341
$ws.Range("D25").Value = 7.591747772672528
$ws.Range("E25").Value = 9.930201710056441
$ws.Range("F25").Value = 37.88976931403817
$ws.Range("I25").Value = 31.00004669632471
$ws.Range("K25").Value = 14.19021587105784
$ws.Range("L25").Value = 10.37231901989365
$ws.Range("N25").Value = 22.11089343246514

Write-Host "Applied 216 cell updates"
